$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 46022.01041666666, 0),
    @(3, 46022.02083333334, 0),
    @(4, 46022.03125, 0),
    @(5, 46022.04166666666, 0),
    @(6, 46022.05208333334, 0),
    @(7, 46022.0625, 0),
    @(8, 46022.07291666666, 0),
    @(9, 46022.08333333334, 0),
    @(10, 46022.09375, 0),
    @(11, 46022.10416666666, 0),
    @(12, 46022.11458333334, 0),
    @(13, 46022.125, 0),
    @(14, 46022.13541666666, 0),
    @(15, 46022.14583333334, 0),
    @(16, 46022.15625, 0),
    @(17, 46022.16666666666, 0),
    @(18, 46022.17708333334, 0.51),
    @(19, 46022.1875, 0),
    @(20, 46022.19791666666, 0),
    @(21, 46022.20833333334, 0),
    @(22, 46022.21875, 0.6840000000000001),
    @(23, 46022.22916666666, 0.681),
    @(24, 46022.23958333334, 0.704),
    @(25, 46022.25, 2.884),
    @(26, 46022.26041666666, 6.92),
    @(27, 46022.27083333334, 10.468),
    @(28, 46022.28125, 14.184),
    @(29, 46022.29166666666, 21.341),
    @(30, 46022.30208333334, 71.91200000000001),
    @(31, 46022.3125, 114.261),
    @(32, 46022.32291666666, 171.181),
    @(33, 46022.33333333334, 237.974),
    @(34, 46022.34375, 438.7),
    @(35, 46022.35416666666, 529.984),
    @(36, 46022.36458333334, 638.1950000000001),
    @(37, 46022.375, 729.162),
    @(38, 46022.38541666666, 912.897),
    @(39, 46022.39583333334, 998.179),
    @(40, 46022.40625, 1074.973),
    @(41, 46022.41666666666, 1141.283),
    @(42, 46022.42708333334, 1220.663),
    @(43, 46022.4375, 1256.403),
    @(44, 46022.44791666666, 1279.585),
    @(45, 46022.45833333334, 1295.029),
    @(46, 46022.46875, 1307.372),
    @(47, 46022.47916666666, 1309.573),
    @(48, 46022.48958333334, 1304.688),
    @(49, 46022.5, 1292.889),
    @(50, 46022.51041666666, 1243.978),
    @(51, 46022.52083333334, 1204.374),
    @(52, 46022.53125, 1153.908),
    @(53, 46022.54166666666, 1094.882),
    @(54, 46022.55208333334, 967.01),
    @(55, 46022.5625, 891.942),
    @(56, 46022.57291666666, 799.79),
    @(57, 46022.58333333334, 707.079),
    @(58, 46022.59375, 518.614),
    @(59, 46022.60416666666, 428.542),
    @(60, 46022.61458333334, 338.513),
    @(61, 46022.625, 259.568),
    @(62, 46022.63541666666, 122.546),
    @(63, 46022.64583333334, 77.411),
    @(64, 46022.65625, 43.632),
    @(65, 46022.66666666666, 26.233),
    @(66, 46022.67708333334, 13.334),
    @(67, 46022.6875, 13.104),
    @(68, 46022.69791666666, 13.424),
    @(69, 46022.70833333334, 13.414),
    @(70, 46022.71875, 7.551),
    @(71, 46022.72916666666, 0),
    @(72, 46022.73958333334, 0),
    @(73, 46022.75, 4.951),
    @(74, 46022.76041666666, 2.05),
    @(75, 46022.77083333334, 1.05),
    @(76, 46022.78125, 0.71),
    @(77, 46022.79166666666, 0),
    @(78, 46022.80208333334, 0),
    @(79, 46022.8125, 0),
    @(80, 46022.82291666666, 0),
    @(81, 46022.83333333334, 0),
    @(82, 46022.84375, 0),
    @(83, 46022.85416666666, 0),
    @(84, 46022.86458333334, 0),
    @(85, 46022.875, 0),
    @(86, 46022.88541666666, 0.61),
    @(87, 46022.89583333334, 0),
    @(88, 46022.90625, 0),
    @(89, 46022.91666666666, 0),
    @(90, 46022.92708333334, 0),
    @(91, 46022.9375, 0),
    @(92, 46022.94791666666, 0),
    @(93, 46022.95833333334, 0),
    @(94, 46022.96875, 0),
    @(95, 46022.97916666666, 0),
    @(96, 46022.98958333334, 0),
    @(97, 46023, 0)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
